$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row count correction: the reworked table only needs 13 data rows
# (3..15) instead of the original 16 (3..18), so drop the 3 trailing rows.
$ws.Rows("16:18").EntireRow.Delete()

# --- Rewrite every data row (3..15) with the corrected / reordered content.
$data = @(
  @(3,  "Combustión Fija",                     "Nafta",             1000,               "MENSUAL", "11/2021"),
  @(4,  "Combustión Móvil",                    "Nafta",             20,                 "MENSUAL", "11/2021"),
  @(5,  "Electricidad adquirida y consumida",  "Electricidad",      200,                "MENSUAL", "12/2021"),
  @(6,  "Electricidad adquirida y consumida",  "Electricidad",      250,                "MENSUAL", "03/2022"),
  @(7,  "Combustión Móvil",                    "GNC",               100,                "MENSUAL", "01/2022"),
  @(8,  "Combustión Fija",                     "Carbón",            260,                "MENSUAL", "07/2022"),
  @(9,  "Combustión Móvil",                    "Gasoil",            150,                "MENSUAL", "06/2022"),
  @(10, "Combustión Fija",                     "Diesel",            1000,               "ANUAL",   "2022"),
  @(11, "Electricidad adquirida y consumida",  "Electricidad",      100,                "MENSUAL", "10/2022"),
  @(12, "Lógistica de productos y servicios",  "peso",              1000,               "MENSUAL", "04/2022"),
  @(13, "Logistica de productos y servicios",  "Distancia",         100,                "MENSUAL", "04/2022"),
  @(14, "Lógistica de productos y servicios",  "categoria",         "materia prima",    "MENSUAL", "04/2022"),
  @(15, "Lógistica de productos y servicios",  "medio_transporte",  "utilitario liviano","MENSUAL", "04/2022")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
  $ws.Range("D$r").Value = $row[4]
  $ws.Range("E$r").Value = $row[5]
}

# --- Highlight row 9 (Combustión Móvil / Gasoil) with underline formatting,
# matching the author's manual emphasis on that record.
$ws.Rows("9").Font.Underline = $true
$ws.Range("A9").WrapText = $true
$ws.Range("B9:D9").WrapText = $true

# --- Column E keeps its special underlined date style only on rows 9 & 12;
# row 13 reverts to the plain (non-underlined) date style.
$ws.Range("E12").Font.Underline = $true
$ws.Range("E13").Font.Underline = $false

# --- View changes: smaller zoom, new active selection.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("A12").Select()
